$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 800.1667
$ws.Range("I9").Value = 760.2
$ws.Range("J9").Value = 1000
$ws.Range("K9").Value = 760.2
$ws.Range("L9").Value = 1000
$ws.Range("M9").Value = -591.2
$ws.Range("N9").Value = -1338
$ws.Range("H40").Value = 8770.076999999999
$ws.Range("J40").Value = 12750.5
$ws.Range("L40").Value = 12750.5
$ws.Range("N40").Value = -13100.5
$ws.Range("H44").Value = 560305
$ws.Range("J44").Value = 9499
$ws.Range("L44").Value = 9499
$ws.Range("N44").Value = -10423
$ws.Range("H62").Value = 7357835
$ws.Range("I62").Value = 11366847
$ws.Range("K62").Value = 11366847
$ws.Range("M62").Value = -11366223
$ws.Range("H65").Value = 7357835
$ws.Range("I65").Value = 11366847
$ws.Range("K65").Value = 56834235
$ws.Range("M65").Value = -56831115
$ws.Range("H106").Value = 2884.7778
$ws.Range("I106").Value = 2884.7778
$ws.Range("K106").Value = 2884.7778
$ws.Range("M106").Value = -2253.7778
$ws.Range("H132").Value = 2647.25
$ws.Range("I132").Value = 2509.0303
$ws.Range("J132").Value = 3298.8572
$ws.Range("K132").Value = 7527.090899999999
$ws.Range("L132").Value = 9896.571599999999
$ws.Range("M132").Value = -4997.090899999999
$ws.Range("N132").Value = -14956.5716
$ws.Range("H135").Value = 715501.6
$ws.Range("J135").Value = 1201
$ws.Range("L135").Value = 10809
$ws.Range("N135").Value = -15879
$ws.Range("H138").Value = 4571.4
$ws.Range("J138").Value = 6095.4033
$ws.Range("L138").Value = 18286.2099
$ws.Range("N138").Value = -28566.2099

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3054.111
$ws.Range("I45").Value = 2248
$ws.Range("K45").Value = 2248
$ws.Range("M45").Value = -1871
$ws.Range("H61").Value = 6686.75
$ws.Range("I61").Value = 3248.75
$ws.Range("K61").Value = 3248.75
$ws.Range("M61").Value = -3036.75
$ws.Range("H110").Value = 126301.98
$ws.Range("I110").Value = 136463.62
$ws.Range("K110").Value = 136463.62
$ws.Range("M110").Value = -134418.62
$ws.Range("H122").Value = 4002.0557
$ws.Range("I122").Value = 2748.762
$ws.Range("J122").Value = 5756.6665
$ws.Range("K122").Value = 8246.286
$ws.Range("L122").Value = 17269.9995
$ws.Range("M122").Value = -5796.286
$ws.Range("N122").Value = -22169.9995
$ws.Range("H132").Value = 4370
$ws.Range("I132").Value = 3618.182
$ws.Range("K132").Value = 10854.546
$ws.Range("M132").Value = -8324.545999999998
$ws.Range("H136").Value = 6686.75
$ws.Range("I136").Value = 3248.75
$ws.Range("K136").Value = 9746.25
$ws.Range("M136").Value = -7196.25

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1664.75
$ws.Range("I86").Value = 1763.7693
$ws.Range("J86").Value = 1235.6666
$ws.Range("K86").Value = 1763.7693
$ws.Range("L86").Value = 1235.6666
$ws.Range("M86").Value = -640.7692999999999
$ws.Range("N86").Value = -3481.6666
$ws.Range("H89").Value = 1664.75
$ws.Range("I89").Value = 1763.7693
$ws.Range("J89").Value = 1235.6666
$ws.Range("K89").Value = 8818.8465
$ws.Range("L89").Value = 6178.333000000001
$ws.Range("M89").Value = -3202.8465
$ws.Range("N89").Value = -17410.333
$ws.Range("H107").Value = 3698.8333
$ws.Range("I107").Value = 3426.6667
$ws.Range("K107").Value = 3426.6667
$ws.Range("M107").Value = -1506.6667

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 88.75
$ws.Range("I6").Value = 46.81818
$ws.Range("J6").Value = 550
$ws.Range("K6").Value = 140.45454
$ws.Range("L6").Value = 1650
$ws.Range("M6").Value = -27.45454000000001
$ws.Range("N6").Value = -1876
$ws.Range("H13").Value = 2141.6667
$ws.Range("I13").Value = 1100
$ws.Range("K13").Value = 3300
$ws.Range("M13").Value = -3132
$ws.Range("H136").Value = 5540.625
$ws.Range("I136").Value = 2865.4
$ws.Range("J136").Value = 9999.333000000001
$ws.Range("K136").Value = 8596.200000000001
$ws.Range("L136").Value = 29997.999
$ws.Range("M136").Value = -3496.200000000001
$ws.Range("N136").Value = -40197.999
$ws.Range("H138").Value = 10001883
$ws.Range("I138").Value = 1854
$ws.Range("K138").Value = 5562
$ws.Range("M138").Value = -422
$ws.Range("H139").Value = 5479.25
$ws.Range("I139").Value = 4160.2383
$ws.Range("K139").Value = 12480.7149
$ws.Range("M139").Value = -7340.714899999999

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 38131.668
$ws.Range("I43").Value = 28990
$ws.Range("J43").Value = 39960
$ws.Range("K43").Value = 28990
$ws.Range("L43").Value = 39960
$ws.Range("M43").Value = -28839
$ws.Range("N43").Value = -40262
$ws.Range("H44").Value = 12440
$ws.Range("J44").Value = 12440
$ws.Range("L44").Value = 12440
$ws.Range("N44").Value = -13632
$ws.Range("H46").Value = 37874.75
$ws.Range("I46").Value = 5000
$ws.Range("K46").Value = 5000
$ws.Range("M46").Value = -4844
$ws.Range("H48").Value = 1666.6666
$ws.Range("I48").Value = 1000
$ws.Range("J48").Value = 2000
$ws.Range("K48").Value = 1000
$ws.Range("L48").Value = 2000
$ws.Range("M48").Value = -515
$ws.Range("N48").Value = -2970
$ws.Range("H53").Value = 0
$ws.Range("J53").Value = 0
$ws.Range("L53").Value = 0
$ws.Range("N53").ClearContents()
$ws.Range("H55").Value = 18714.285
$ws.Range("I55").Value = 20000
$ws.Range("J55").Value = 18200
$ws.Range("K55").Value = 20000
$ws.Range("L55").Value = 18200
$ws.Range("M55").Value = -19673
$ws.Range("N55").Value = -18854
$ws.Range("H58").Value = 30046
$ws.Range("I58").Value = 0
$ws.Range("K58").Value = 0
$ws.Range("M58").ClearContents()
$ws.Range("H59").Value = 24109.666
$ws.Range("I59").Value = 0
$ws.Range("J59").Value = 24109.666
$ws.Range("K59").Value = 0
$ws.Range("L59").Value = 24109.666
$ws.Range("M59").ClearContents()
$ws.Range("N59").Value = -25275.666
$ws.Range("H74").Value = 235065
$ws.Range("J74").Value = 235065
$ws.Range("L74").Value = 235065
$ws.Range("N74").Value = -236937
$ws.Range("H75").Value = 95086.664
$ws.Range("J75").Value = 95086.664
$ws.Range("L75").Value = 95086.664
$ws.Range("N75").Value = -96834.664
$ws.Range("H77").Value = 235065
$ws.Range("J77").Value = 235065
$ws.Range("L77").Value = 705195
$ws.Range("N77").Value = -714555
$ws.Range("H78").Value = 95086.664
$ws.Range("J78").Value = 95086.664
$ws.Range("L78").Value = 285259.992
$ws.Range("N78").Value = -293995.992
$ws.Range("H80").Value = 1671421.8
$ws.Range("I80").Value = 1115839.9
$ws.Range("J80").Value = 3338167.2
$ws.Range("K80").Value = 1115839.9
$ws.Range("L80").Value = 3338167.2
$ws.Range("M80").Value = -1114841.9
$ws.Range("N80").Value = -3340163.2
$ws.Range("H82").Value = 0
$ws.Range("J82").Value = 0
$ws.Range("L82").Value = 0
$ws.Range("N82").ClearContents()
$ws.Range("H83").Value = 1671421.8
$ws.Range("I83").Value = 1115839.9
$ws.Range("J83").Value = 3338167.2
$ws.Range("K83").Value = 5579199.5
$ws.Range("L83").Value = 16690836
$ws.Range("M83").Value = -5574207.5
$ws.Range("N83").Value = -16700820
$ws.Range("H85").Value = 0
$ws.Range("J85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("N85").ClearContents()
$ws.Range("H132").Value = 92525.586
$ws.Range("I132").Value = 11458.8
$ws.Range("J132").Value = 150430.42
$ws.Range("K132").Value = 34376.39999999999
$ws.Range("L132").Value = 451291.26
$ws.Range("M132").Value = -31846.39999999999
$ws.Range("N132").Value = -456351.26

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 5520.8887
$ws.Range("I61").Value = 5138.4
$ws.Range("J61").Value = 5999
$ws.Range("K61").Value = 5138.4
$ws.Range("L61").Value = 5999
$ws.Range("M61").Value = -4936.4
$ws.Range("N61").Value = -6403
$ws.Range("H100").Value = 3748.1667
$ws.Range("I100").Value = 3622.25
$ws.Range("K100").Value = 3622.25
$ws.Range("M100").Value = -3081.25
$ws.Range("H113").Value = 5520.8887
$ws.Range("I113").Value = 5138.4
$ws.Range("J113").Value = 5999
$ws.Range("K113").Value = 5138.4
$ws.Range("L113").Value = 5999
$ws.Range("M113").Value = -2968.4
$ws.Range("N113").Value = -10339

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H74").Value = 6926.0835
$ws.Range("J74").Value = 7636
$ws.Range("L74").Value = 7636
$ws.Range("N74").Value = -9508
$ws.Range("H77").Value = 6926.0835
$ws.Range("J77").Value = 7636
$ws.Range("L77").Value = 22908
$ws.Range("N77").Value = -32268
$ws.Range("H100").Value = 961.05884
$ws.Range("I100").Value = 1038.7858
$ws.Range("K100").Value = 2077.5716
$ws.Range("M100").Value = -1536.5716
$ws.Range("H136").Value = 10468006
$ws.Range("I136").Value = 15608943
$ws.Range("J136").Value = 186132.9
$ws.Range("K136").Value = 46826829
$ws.Range("L136").Value = 558398.7
$ws.Range("M136").Value = -46824279
$ws.Range("N136").Value = -563498.7
